# "Tried to implement Penality Reward System (unfinished)"
#
# This script mutates the "Daily PO" sheet (rows 9-16 get their PO records
# rotated/reshuffled - Window start/end travel with the moved record, but the
# trailing Index column stays put) and drops the now-obsolete last PO row
# (row 19, PO 86F3MF6J). The downstream "Merged (Optional)" rollup sheet,
# "PO Volume Insights" aggregates and the "PO Prediction" forecast are all
# updated to stay consistent with the new Daily_PO_Qty numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Daily PO" sheet: reshuffle rows 9,10,11,12,14,15,16, then drop row 19
# ---------------------------------------------------------------------
$wsPO = $wb.Worksheets.Item("Daily PO")

# Row 9  (was 1BZWCQBV/LGB8, now takes row10's PO data; T/U unaffected here)
$wsPO.Range("A9").Value = "5I9ZRHZV"
$wsPO.Range("M9").Value = 140
$wsPO.Range("N9").Value = 140
$wsPO.Range("P9").Value = 140
$wsPO.Range("S9").Value = "SBD1"
$wsPO.Range("Y9").Value = 16100
$wsPO.Range("Z9").Value = 16100
$wsPO.Range("AA9").Value = 16100

# Row 10 (was 5I9ZRHZV/SBD1, now takes row11's PO data incl. window dates)
$wsPO.Range("A10").Value = "74FIAHPZ"
$wsPO.Range("M10").Value = 20
$wsPO.Range("N10").Value = 20
$wsPO.Range("P10").Value = 20
$wsPO.Range("S10").Value = "MDW2"
$wsPO.Range("T10").Value = 45041
$wsPO.Range("U10").Value = 45043
$wsPO.Range("Y10").Value = 2300
$wsPO.Range("Z10").Value = 2300
$wsPO.Range("AA10").Value = 2300

# Row 11 (was 74FIAHPZ/MDW2, now takes row12's PO data)
$wsPO.Range("A11").Value = "7EAUKLWC"
$wsPO.Range("M11").Value = 10
$wsPO.Range("N11").Value = 10
$wsPO.Range("P11").Value = 10
$wsPO.Range("S11").Value = "FWA4"
$wsPO.Range("Y11").Value = 1150
$wsPO.Range("Z11").Value = 1150
$wsPO.Range("AA11").Value = 1150

# Row 12 (was 7EAUKLWC/FWA4, now takes row16's PO data, incl. cancelled qty)
$wsPO.Range("A12").Value = "7UC1XMLE"
$wsPO.Range("M12").Value = 80
$wsPO.Range("N12").Value = 80
$wsPO.Range("P12").Value = 80
$wsPO.Range("Q12").Value = 2
$wsPO.Range("S12").Value = "SCK4"
$wsPO.Range("Y12").Value = 9200
$wsPO.Range("Z12").Value = 9200
$wsPO.Range("AA12").Value = 9200
$wsPO.Range("AB12").Value = 230

# Row 13 is untouched (7MM8EEPT stays as-is)

# Row 14 (was 7W6SJQGT/LAX9, now takes row15's PO data incl. window dates)
$wsPO.Range("A14").Value = "7Y6ZLXLI"
$wsPO.Range("M14").Value = 30
$wsPO.Range("N14").Value = 30
$wsPO.Range("P14").Value = 30
$wsPO.Range("T14").Value = 45009
$wsPO.Range("U14").Value = 45013
$wsPO.Range("Y14").Value = 3450
$wsPO.Range("Z14").Value = 3450
$wsPO.Range("AA14").Value = 3450

# Row 15 (was 7Y6ZLXLI/LAX9, now takes row9's original PO data incl. dates)
$wsPO.Range("A15").Value = "1BZWCQBV"
$wsPO.Range("M15").Value = 40
$wsPO.Range("N15").Value = 40
$wsPO.Range("P15").Value = 40
$wsPO.Range("S15").Value = "LGB8"
$wsPO.Range("T15").Value = 44965
$wsPO.Range("U15").Value = 44972
$wsPO.Range("Y15").Value = 4600
$wsPO.Range("Z15").Value = 4600
$wsPO.Range("AA15").Value = 4600

# Row 16 (was 7UC1XMLE/SCK4, now takes row14's original PO data)
$wsPO.Range("A16").Value = "7W6SJQGT"
$wsPO.Range("M16").Value = 40
$wsPO.Range("N16").Value = 40
$wsPO.Range("P16").Value = 40
$wsPO.Range("Q16").Value = 0
$wsPO.Range("S16").Value = "LAX9"
$wsPO.Range("Y16").Value = 4600
$wsPO.Range("Z16").Value = 4600
$wsPO.Range("AA16").Value = 4600
$wsPO.Range("AB16").Value = 0

# Row 19 (PO 86F3MF6J, cancelled/out-of-stock duplicate of row18's date+qty)
# is removed entirely; everything below it (nothing, it was last) shifts up.
$wsPO.Rows.Item(19).Delete()

# ---------------------------------------------------------------------
# 2) "Merged (Optional)" sheet: mirror the new Daily_PO_Qty values (col C,
#    rows 11-18 line up with "Daily PO" rows 9-16) and drop the row that
#    corresponded to the deleted PO (the duplicate 45133 / qty 200 entry).
# ---------------------------------------------------------------------
$wsMerged = $wb.Worksheets.Item("Merged (Optional)")

$wsMerged.Range("C11").Value = 140
$wsMerged.Range("C12").Value = 20
$wsMerged.Range("C13").Value = 10
$wsMerged.Range("C14").Value = 80
$wsMerged.Range("C15").Value = 180
$wsMerged.Range("C16").Value = 30
$wsMerged.Range("C17").Value = 40
$wsMerged.Range("C18").Value = 40

$wsMerged.Rows.Item(30).Delete()

# ---------------------------------------------------------------------
# 3) "PO Volume Insights" sheet: totals/averages over the updated
#    Daily_PO_Qty column (max/min are unchanged: still 200 / 10).
# ---------------------------------------------------------------------
$wsInsights = $wb.Worksheets.Item("PO Volume Insights")

$wsInsights.Range("A2").Value = 1030
$wsInsights.Range("B2").Value = 60.58823529411764

# ---------------------------------------------------------------------
# 4) "PO Prediction" sheet: next-period forecast recomputed from the
#    updated Daily_PO_Qty series.
# ---------------------------------------------------------------------
$wsPrediction = $wb.Worksheets.Item("PO Prediction")

$wsPrediction.Range("A2").Value = 94.77941176470588
